# Slide 15 ("blank" content slide) still carried two empty, unused
# placeholder shapes (an orphaned Title and Content placeholder) left over
# from the layout. Remove them so the slide's shape tree is empty, matching
# the cleaned-up deck.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $s.Shapes.Item($i).Delete()
}
